$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Insert a new row above the current row 3 (the sub-header row with
# FC / Default / etc.) so it shifts down to become row 4, leaving a
# fresh blank row 3 for the new "financial cost component" headers.
# ------------------------------------------------------------------
$ws.Rows(3).Insert()

# ------------------------------------------------------------------
# Row 2 - the "Financial Cost" banner becomes "Components of FC" and
# the merge shrinks from B2:G2 to C2:G2 (column B goes back to a
# plain, unmerged cell like A2/I2).
# ------------------------------------------------------------------
$ws.Range("B2:G2").UnMerge()
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = "Components of FC"
$ws.Range("C2:G2").Merge()

# B2 / I2 drop their centering and pick up a thin top border (same
# look as A2).
$ws.Range("B2").HorizontalAlignment = -4131
$ws.Range("B2").Borders.Item(8).LineStyle = 1
$ws.Range("B2").Borders.Item(8).Weight = 2
$ws.Range("B2").Borders.Item(9).LineStyle = -4142

$ws.Range("I2").HorizontalAlignment = -4131

# C2:G2 keep their centering but now get a thin top AND thin bottom
# border (was: no top border / double bottom border).
$rng = $ws.Range("C2:G2")
$rng.Borders.Item(8).LineStyle = 1
$rng.Borders.Item(8).Weight = 2
$rng.Borders.Item(9).LineStyle = 1
$rng.Borders.Item(9).Weight = 2

# ------------------------------------------------------------------
# New row 3: column sub-headers for the financial-cost decomposition.
# ------------------------------------------------------------------
$ws.Range("B3").Value = "FC"
$ws.Range("C3").Value = "Interest pymnt"
$ws.Range("D3").Value = "Fee pymnt"
$ws.Range("E3").Value = "Principal pymnt"
$ws.Range("F3").Value = "Lost pawn value"
$ws.Range("G3").Value = "Default"
$ws.Range("I3").Value = "APR"

$hdrRng = $ws.Range("B3:I3")
$hdrRng.HorizontalAlignment = -4108

$ws.Range("B3").Borders.Item(9).LineStyle = 1
$ws.Range("B3").Borders.Item(9).Weight = 2
$ws.Range("I3").Borders.Item(9).LineStyle = 1
$ws.Range("I3").Borders.Item(9).Weight = 2

# ------------------------------------------------------------------
# Row 4 (old row 3, pushed down by the insert): drop the FC / Lost
# pawn value / Default headers that moved up to row 3, and replace
# them with the LaTeX formula labels for each column.
# ------------------------------------------------------------------
$ws.Range("B4").ClearContents()
$ws.Range("C4").Value = "\sum_t P^i_{it}"
$ws.Range("D4").Value = "\sum_t P^f_{it}"
$ws.Range("E4").Value = "\mathds{1}(\text{Def}_i)}\times\sum_t P^c_{it}"
$ws.Range("F4").Value = "\mathds{1}(\text{Def}_i)}\times \text{Appr. Val.}_i"
$ws.Range("G4").Value = "\mathds{1}(\text{Def}_i)"

# H4 / I4 lose their centered-but-borderless look; I4 also loses the
# stray "APR" value that the row-insert carried down from the old I3
# (that label now lives in the new I3 created above).
$ws.Range("H4").HorizontalAlignment = -4131
$ws.Range("I4").HorizontalAlignment = -4131
$ws.Range("I4").ClearContents()
